# 03_Preparation_Answer_Key_C.docx edit
# - Append " (Questions 1-5)" to the bold instructional paragraph, plus a
#   trailing bold space run.
# - Delete the "1 / Mode / The most frequently occurring value" table row.
# - Strip the " Mode -NN.NNNNNN" suffix from each Company A-E solution cell.

$d = $word.ActiveDocument

# 1) Extend the bold note paragraph and add a trailing space run.
$notePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -match "^Please note that the steps show rounded numbers") {
        $notePara = $candidate
        break
    }
}

$noteRange = $notePara.Range
$noteRange.End = $noteRange.End - 1
$noteRange.Text = "Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding. (Questions 1-5)"

$tail = $d.Range($notePara.Range.End - 1, $notePara.Range.End - 1)
$tail.InsertAfter(" ")

# 2) Remove the Problem 1 "Mode" row from the solutions table.
$table = $d.Tables.Item(1)
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Item(2).Range.Text -match "^Mode") {
        $row.Delete()
        break
    }
}

# 3) Drop the trailing " Mode -NN.NNNNNN" readout from each company row.
$replacements = @(
    @{company = "Company A"; new = "Mean: 21.276 Median: 13.433"},
    @{company = "Company B"; new = "Mean: 33.482 Median: 20.838"},
    @{company = "Company C"; new = "Mean: 41.122 Median: 25.558"},
    @{company = "Company D"; new = "Mean: 0.706 Median: 1.892"},
    @{company = "Company E"; new = "Mean: -1.084 Median: -3.796"}
)

for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    $nameCell = $row.Cells.Item(2).Range.Text
    foreach ($rep in $replacements) {
        if ($nameCell -match $rep.company) {
            $cellRange = $row.Cells.Item(3).Range
            $cellRange.End = $cellRange.End - 2
            $cellRange.Text = $rep.new
        }
    }
}
